# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8 and 9 swapped their Coin/Link content (Cardano <-> OKB)
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"

# Price (column D) and Volume(1h) (column E) refresh for every data row
$updates = @(
    @{ Row = 2; D = "24.821.34"; E = "  +2.01%  " },
    @{ Row = 3; D = "1.716.46"; E = "  +1.85%  " },
    @{ Row = 4; D = "1.011"; E = "  +1.10%  " },
    @{ Row = 5; D = "312.00"; E = "  +1.34%  " },
    @{ Row = 6; D = "1.006"; E = "  +0.94%  " },
    @{ Row = 7; D = "0.3779"; E = "  +1.39%  " },
    @{ Row = 8; D = "49.66"; E = "  +2.77%  " },
    @{ Row = 9; D = "0.3490"; E = "  +1.58%  " },
    @{ Row = 10; D = "1.192"; E = "  +0.56%  " },
    @{ Row = 11; D = "0.07473"; E = "  +2.29%  " },
    @{ Row = 12; D = "1.009"; E = "  +1.25%  " },
    @{ Row = 13; D = "6.346"; E = "  +4.02%  " },
    @{ Row = 14; D = "20.84"; E = "  +1.08%  " },
    @{ Row = 15; D = "6.982"; E = "  +3.30%  " },
    @{ Row = 16; D = "1.722.96"; E = "  +2.97%  " },
    @{ Row = 17; D = "0.00001127"; E = "  +1.70%  " },
    @{ Row = 18; D = "1.008"; E = "  +0.97%  " },
    @{ Row = 19; D = "0.06679"; E = "  -0.64%  " },
    @{ Row = 20; D = "83.97"; E = "  +3.05%  " },
    @{ Row = 21; D = "17.26"; E = "  +4.71%  " },
    @{ Row = 22; D = "6.378"; E = "  +3.98%  " },
    @{ Row = 23; D = "13.52"; E = "  +12.25%  " },
    @{ Row = 24; D = "24.835.58"; E = "  +2.42%  " },
    @{ Row = 25; D = "2.448"; E = "  +1.87%  " },
    @{ Row = 26; D = "2.804"; E = "  +4.54%  " },
    @{ Row = 27; D = "20.44"; E = "  +4.60%  " },
    @{ Row = 28; D = "150.80"; E = "  -0.98%  " },
    @{ Row = 29; D = "1.913.26"; E = "  +2.96%  " },
    @{ Row = 30; D = "131.76"; E = "  +3.67%  " },
    @{ Row = 31; D = "1.173"; E = "  +18.39%  " },
    @{ Row = 32; D = "6.830"; E = "  +5.67%  " },
    @{ Row = 33; D = "4.236"; E = "  +4.94%  " },
    @{ Row = 34; D = "1.798"; E = "  +3.17%  " },
    @{ Row = 35; D = $null; E = "  +9.81%  " },
    @{ Row = 36; D = "0.08792"; E = "  +3.94%  " },
    @{ Row = 37; D = "5.641"; E = "  +4.69%  " },
    @{ Row = 38; D = "0.02438"; E = "  +4.17%  " },
    @{ Row = 39; D = "0.06538"; E = "  +1.75%  " },
    @{ Row = 40; D = $null; E = "  +1.19%  " },
    @{ Row = 41; D = "0.2206"; E = "  +4.36%  " },
    @{ Row = 42; D = $null; E = "  -1.43%  " },
    @{ Row = 43; D = "0.6450"; E = "  +4.77%  " },
    @{ Row = 44; D = "1.007"; E = "  +1.01%  " },
    @{ Row = 45; D = "13.88"; E = "  +5.26%  " },
    @{ Row = 46; D = "0.6141"; E = "  +2.98%  " },
    @{ Row = 47; D = "3.848"; E = "  +1.26%  " },
    @{ Row = 48; D = "2.148"; E = "  +6.17%  " },
    @{ Row = 49; D = "129.32"; E = "  +1.07%  " },
    @{ Row = 50; D = "0.07275"; E = "  +1.53%  " },
    @{ Row = 51; D = "79.84"; E = "  +4.00%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Prefix with an apostrophe so Excel keeps these price strings as text
        # (many look like numbers, e.g. "1.011") instead of coercing them to
        # doubles, then strip the resulting quote-prefix formatting so the
        # cell style stays untouched, matching the original workbook.
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.Value = "'" + $u.D
        $cell.ClearFormats()
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
